# Revert the presentation's applied design theme from the custom
# "Integral" colour scheme back to the stock Office "Office Theme"
# colour scheme. The deck has a single slide master/design, whose theme
# part (ppt/theme/theme2.xml) is shared by every slide, so a single
# update through the Design's colour scheme recolors the whole deck.
#
# Target colours (standard Office 2019+ "Office Theme" clrScheme),
# expressed as VBA-style BGR-packed RGB() long values:
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
